$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width change from 16.42578125 to 15.42578125 (raw OOXML "width" units).
# Excel's ColumnWidth COM property is character-width based and gets rounded to
# the nearest pixel when converted back to the raw stored width, so we pick the
# ColumnWidth value whose stored width lands closest to the target 15.42578125.
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666

$ws.Range("A1").Value = -0.14907911348512926
$ws.Range("B1").Value = 0.14834444936677471
$ws.Range("A2").Value = -0.02576572683877032
$ws.Range("B2").Value = 0.024446879909712393
$ws.Range("A3").Value = 0.078490584482821646
$ws.Range("B3").Value = -0.079184709825096178
$ws.Range("A4").Value = -0.20480893633784092
$ws.Range("B4").Value = 0.20328569544572872
$ws.Range("A5").Value = -0.19728569585396372
$ws.Range("B5").Value = 0.19417992706687492
$ws.Range("A6").Value = -0.093253386831670948
$ws.Range("B6").Value = 0.093125183279025858
$ws.Range("A7").Value = -0.073125183778573799
$ws.Range("B7").Value = 0.072829357918768878
$ws.Range("A8").Value = -0.052829358422760819
$ws.Range("B8").Value = 0.05260853092815676
$ws.Range("A9").Value = -0.046608531365825989
$ws.Range("B9").Value = 0.046432167603966512
$ws.Range("A10").Value = -0.040432168046464767
$ws.Range("B10").Value = 0.040410617554307748
$ws.Range("A11").Value = -0.035910617989429738
$ws.Range("B11").Value = 0.035875320851776138
$ws.Range("A12").Value = -0.029875321295653734
$ws.Range("B12").Value = 0.029776684383030716
$ws.Range("A13").Value = -0.023776684831364747
$ws.Range("B13").Value = 0.023754575152040047
$ws.Range("A14").Value = -0.011754575632372699
$ws.Range("B14").Value = 0.011749235638448319
$ws.Range("A15").Value = -0.02105389400485258
$ws.Range("B15").Value = 0.021028120499693159
$ws.Range("A16").Value = -0.01502812095115913
$ws.Range("B16").Value = 0.015004549754051855
$ws.Range("A17").Value = -0.0090045502074254102
$ws.Range("B17").Value = 0.0089999995308014391
$ws.Range("A18").Value = -0.076485574154233404
$ws.Range("B18").Value = 0.076347599939040833
$ws.Range("A19").Value = -0.067347600357219317
$ws.Range("B19").Value = 0.066244097979569982
$ws.Range("A20").Value = -0.071929554227580184
$ws.Range("B20").Value = 0.07168307808406027
$ws.Range("A21").Value = -0.0090044233321906475
$ws.Range("B21").Value = 0.0089999995692084944
$ws.Range("A22").Value = -0.093937464811579829
$ws.Range("B22").Value = 0.093627675110530717
$ws.Range("A23").Value = -0.084627675532678026
$ws.Range("B23").Value = 0.084125470669481395
$ws.Range("A24").Value = -0.042125471268317227
$ws.Range("B24").Value = 0.041999999398070642
$ws.Range("A25").Value = -0.094926540661024461
$ws.Range("B25").Value = 0.094682068940464603
$ws.Range("A26").Value = -0.088682069367717276
$ws.Range("B26").Value = 0.088367027369208984
$ws.Range("A27").Value = -0.082367027798712744
$ws.Range("B27").Value = 0.081289171513187597
$ws.Range("A28").Value = -0.075289171952226397
$ws.Range("B28").Value = 0.074542339161467908
$ws.Range("A29").Value = -0.062542339637863265
$ws.Range("B29").Value = 0.062171991708046548
$ws.Range("A30").Value = -0.0421719922285404
$ws.Range("B30").Value = 0.042020196506194374
$ws.Range("A31").Value = -0.027020197005354518
$ws.Range("B31").Value = 0.027000848573456082
$ws.Range("A32").Value = -0.0060008491041090295
$ws.Range("B32").Value = 0.0059999995462440836
